# Update "江西-漫展信息" workbook: refresh a handful of "want-to-go" counts
# and insert a newly-announced event (宜春·原x穹x崩only) as row 16 in the
# two sheets that carry the full event table (展览 / 全部类型).

function Update-ExpoSheet($ws) {
    # --- 1. Insert a new row at position 16, pushing the old row 16 (and
    #        everything below it) down by one. This also grows the used
    #        range to A1:I28 automatically. ---
    $ws.Rows.Item(16).Insert()

    # The freshly inserted row comes back blank/unstyled except for a
    # default style on A16, so rebuild the same look used by every other
    # row-number cell in column A (bold, centered, thin box border).
    $ws.Cells.Item(16, 1).Value = 15
    $ws.Cells.Item(16, 1).Font.Bold = $true
    $ws.Cells.Item(16, 1).HorizontalAlignment = -4108
    $ws.Cells.Item(16, 1).VerticalAlignment = -4160
    $ws.Cells.Item(16, 1).Borders.LineStyle = 1

    # --- 2. Fill in the new event's details (B16:I16). ---
    $ws.Cells.Item(16, 2).Value = "2024-04-05"
    $ws.Cells.Item(16, 3).Value = "宜春·原x穹x崩only"
    $ws.Cells.Item(16, 4).Value = "市府北路10号  红林大酒店"
    $ws.Cells.Item(16, 5).Value = "2024.04.05 10:00-04.05 17:00"
    $ws.Cells.Item(16, 6).Value = 1
    $ws.Cells.Item(16, 7).Value = 35
    $ws.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83073"
    $ws.Cells.Item(16, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/l8vN3pmn1710918987797.jpeg"

    # --- 3. Refresh "want to go" counts (column F) that ticked up since
    #        the last export. Rows below the insertion point are addressed
    #        at their *new* (post-insert) row numbers. ---
    $ws.Cells.Item(3, 6).Value = 3040
    $ws.Cells.Item(7, 6).Value = 1660
    $ws.Cells.Item(9, 6).Value = 84
    $ws.Cells.Item(10, 6).Value = 33
    $ws.Cells.Item(11, 6).Value = 1362
    $ws.Cells.Item(13, 6).Value = 498
    $ws.Cells.Item(15, 6).Value = 26
    $ws.Cells.Item(17, 6).Value = 73
    $ws.Cells.Item(18, 6).Value = 62
    $ws.Cells.Item(19, 6).Value = 123
    $ws.Cells.Item(22, 6).Value = 3168
    $ws.Cells.Item(23, 6).Value = 388
    $ws.Cells.Item(24, 6).Value = 122
    $ws.Cells.Item(25, 6).Value = 204
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Update-ExpoSheet $ws
}
